$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format before writing values so that
# numeric-looking strings (e.g. "1.00", "0.0000101") keep their exact
# textual representation instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.783.96'
$ws.Range("E2").Value = '  -2.82%  '
$ws.Range("D3").Value = '2.491.17'
$ws.Range("E3").Value = '  -5.06%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '554.01'
$ws.Range("E5").Value = '  -3.66%  '
$ws.Range("D6").Value = '147.68'
$ws.Range("E6").Value = '  -4.42%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.605'
$ws.Range("E8").Value = '  -3.01%  '
$ws.Range("D9").Value = '2.491.42'
$ws.Range("E9").Value = '  -4.94%  '
$ws.Range("E10").Value = '  -7.45%  '
$ws.Range("D11").Value = '5.45'
$ws.Range("E11").Value = '  -5.94%  '
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("E13").Value = '  -5.14%  '
$ws.Range("D14").Value = '26.26'
$ws.Range("E14").Value = '  -6.99%  '
$ws.Range("D15").Value = '2.942.84'
$ws.Range("E15").Value = '  -4.87%  '
$ws.Range("E16").Value = '  -8.10%  '
$ws.Range("D17").Value = '61.706.16'
$ws.Range("E17").Value = '  -2.76%  '
$ws.Range("D18").Value = '2.490.43'
$ws.Range("E18").Value = '  -4.32%  '
$ws.Range("D19").Value = '11.21'
$ws.Range("E19").Value = '  -7.14%  '
$ws.Range("D20").Value = '7.01'
$ws.Range("E20").Value = '  -6.90%  '
$ws.Range("D21").Value = '4.22'
$ws.Range("E21").Value = '  -6.41%  '
$ws.Range("D22").Value = '323.34'
$ws.Range("E22").Value = '  -5.94%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -3.14%  '
$ws.Range("D25").Value = '64.15'
$ws.Range("E25").Value = '  -5.38%  '
$ws.Range("D26").Value = '0.0000101'
$ws.Range("E26").Value = '  -5.67%  '
$ws.Range("D27").Value = '2.609.23'
$ws.Range("E27").Value = '  -4.80%  '
$ws.Range("D28").Value = '1.54'
$ws.Range("E28").Value = '  -3.49%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '8.44'
$ws.Range("E30").Value = '  -7.96%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '537.45'
$ws.Range("E31").Value = '  -9.64%  '
$ws.Range("D32").Value = '7.61'
$ws.Range("E32").Value = '  -3.69%  '
$ws.Range("D33").Value = '0.153'
$ws.Range("E33").Value = '  -4.95%  '
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  -6.77%  '
$ws.Range("E35").Value = '  -7.39%  '
$ws.Range("D36").Value = '6.02'
$ws.Range("E36").Value = '  -8.24%  '
$ws.Range("D37").Value = '4.92'
$ws.Range("E37").Value = '  -8.20%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("E39").Value = '  -3.59%  '
$ws.Range("E40").Value = '  -5.44%  '
$ws.Range("D41").Value = '148.61'
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("E42").Value = '  -7.98%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '40.38'
$ws.Range("E44").Value = '  -3.06%  '
$ws.Range("E45").Value = '  -5.26%  '
$ws.Range("D46").Value = '149.16'
$ws.Range("E46").Value = '  -5.75%  '
$ws.Range("E47").Value = '  -6.59%  '
$ws.Range("D48").Value = '21.02'
$ws.Range("E48").Value = '  -13.06%  '
$ws.Range("E49").Value = '  -8.28%  '
$ws.Range("D50").Value = '0.600'
$ws.Range("E50").Value = '  -4.47%  '
$ws.Range("D51").Value = '0.0949'
$ws.Range("E51").Value = '  -4.89%  '
